# Apply the cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "0.999", "28.80", "1.00") must be pre-formatted as Text so the literal
# string (including insignificant trailing zeros) survives the round trip.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Write the updated values
$ws.Range('D2').Value = '59.259.21'
$ws.Range('E2').Value = '  +4.74%  '
$ws.Range('D3').Value = '3.348.60'
$ws.Range('E3').Value = '  +2.70%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '410.25'
$ws.Range('E5').Value = '  +2.97%  '
$ws.Range('D6').Value = '112.71'
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('D7').Value = '0.586'
$ws.Range('E7').Value = '  +4.75%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.638'
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('D10').Value = '40.35'
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').Value = '0.100'
$ws.Range('E11').Value = '  +3.94%  '
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('D13').Value = '3.874.85'
$ws.Range('E13').Value = '  +2.84%  '
$ws.Range('D14').Value = '8.55'
$ws.Range('E14').Value = '  +4.86%  '
$ws.Range('D15').Value = '19.44'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '3.351.45'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '59.048.68'
$ws.Range('E18').Value = '  +4.53%  '
$ws.Range('D19').Value = '10.83'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('D20').Value = '3.37'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').Value = '0.0000113'
$ws.Range('E21').Value = '  +7.62%  '
$ws.Range('D22').Value = '13.19'
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('D23').Value = '305.38'
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('D24').Value = '75.78'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('D25').Value = '3.19'
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').Value = '28.80'
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('D28').Value = '7.92'
$ws.Range('E28').Value = '  -2.77%  '
$ws.Range('D29').Value = '7.60'
$ws.Range('E29').Value = '  +3.15%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.174'
$ws.Range('E30').Value = '  +2.26%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.118'
$ws.Range('E31').Value = '  +6.18%  '
$ws.Range('D32').Value = '11.60'
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('D34').Value = '40.78'
$ws.Range('E34').Value = '  +9.49%  '
$ws.Range('D35').Value = '0.0525'
$ws.Range('E35').Value = '  +7.49%  '
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('D37').Value = '52.06'
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '3.09'
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('D40').Value = '3.42'
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('D41').Value = '137.63'
$ws.Range('E41').Value = '  +3.17%  '
$ws.Range('D42').Value = '0.123'
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('D43').Value = '1.94'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').Value = '4.02'
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('D45').Value = '17.04'
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('D46').Value = '0.281'
$ws.Range('E46').Value = '  -2.10%  '
$ws.Range('D47').Value = '2.26'
$ws.Range('E47').Value = '  +7.96%  '
$ws.Range('D48').Value = '22.34'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').Value = '2.208.13'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('D50').Value = '2.40'
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').Value = '6.37'
$ws.Range('E51').Value = '  +5.96%  '

Write-Host "cryptos list updated"
